$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Ringgold"
$ws.Range("C7").Value = "Ringgold identifier for organisations in the publishing industry supply chain"
